$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Remove the duplicated "Contact / No display for ContactDetail" row (row 11),
# shifting all subsequent rows up by one.
$ws.Rows.Item(11).Delete()

# Version bump
$ws.Range("B3").Value = "6.0.0"

# Date refresh
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value
$ws.Range("B9").Value = "Alvearie Team"

# Former "Contact" row (now row 10) becomes "Jurisdiction"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
